$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 ---
# B10: value changed, format unchanged (0.0000%)
$ws.Range("B10").Value = 0.75510200000000005
# C10: value -> 1, format changes from 0.0000% to 0% (built-in percent)
$ws.Range("C10").Value = 1
$ws.Range("C10").NumberFormat = "0%"

# --- Row 11 ---
# B11: value changed only
$ws.Range("B11").Value = 0.68965500000000002

# --- Row 12 ---
# B12: value unchanged, format changes from 0.0000% to 0.000%
$ws.Range("B12").NumberFormat = "0.000%"

# --- Row 14 ---
# B14: value changed only
$ws.Range("B14").Value = 0.63636400000000004

# --- Row 15 ---
# C15: value changes, format changes from 0.0000% to a new 0.0% custom format
$ws.Range("C15").Value = 0.996
$ws.Range("C15").NumberFormat = "0.0%"

# --- Row 16 ---
# B16: value changed only
$ws.Range("B16").Value = 0.90163899999999997

# --- Row 17 ---
# C17: value changed only
$ws.Range("C17").Value = 0.965665

# --- Row 18 ---
# B18: value changes, format changes from 0% to 0.0000%
$ws.Range("B18").Value = 0.64285700000000001
$ws.Range("B18").NumberFormat = "0.0000%"

# --- Row 19 ---
# B19: value changes, format changes from 0.000% to 0.0000%
$ws.Range("B19").Value = 0.70422499999999999
$ws.Range("B19").NumberFormat = "0.0000%"
